$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newQuery = "MATCH (p:participant)-->(s:study)`nOPTIONAL MATCH (samp:sample)-->(p)`nOPTIONAL MATCH (p)<--(diag:diagnosis)`nOPTIONAL MATCH (samp)<--(f:file)`nOPTIONAL MATCH (f)<--(g:genomic_info)`nWITH s, p, samp, f, g, diag`nWHERE diag.primary_diagnosis in ['Adrenal Cortical Tumor']`nWITH p`nOPTIONAL MATCH (p)-->(s:study)`nOPTIONAL MATCH (samp:sample)-->(p)`nWITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp`nRETURN `ncoalesce(p.participant_id,'') as ``Participant ID```,`ncoalesce(s.study_name, '') as ``Study Name```,`ncoalesce(s.phs_accession,'') as ``Accession```,`ncoalesce(p.gender,'') as ``Gender```,`ncoalesce(apoc.text.join(samp, ','), '') as ``Samples```nORDER BY p.participant_id`nLIMIT 100"

$ws.Range("B2").Value = $newQuery
$ws.Rows.Item(2).RowHeight = 157.5

$ws.Range("C2").Select()
